$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7150.952
$ws.Range("J40").Value = 7617.8
$ws.Range("L40").Value = 7617.8
$ws.Range("N40").Value = -7967.8
$ws.Range("H41").Value = 20836362
$ws.Range("I41").Value = 515.5833
$ws.Range("K41").Value = 515.5833
$ws.Range("M41").Value = -75.58330000000001
$ws.Range("H49").Value = 264.5
$ws.Range("J49").Value = 129
$ws.Range("L49").Value = 387
$ws.Range("N49").Value = -659
$ws.Range("H96").Value = 626.5
$ws.Range("I96").Value = 513.1111
$ws.Range("K96").Value = 1539.3333
$ws.Range("M96").Value = -166.3332999999998
$ws.Range("H97").Value = 2779.6155
$ws.Range("J97").Value = 2928
$ws.Range("L97").Value = 8784
$ws.Range("N97").Value = -9776
$ws.Range("H111").Value = 62989.766
$ws.Range("I111").Value = 115382.664
$ws.Range("K111").Value = 346147.992
$ws.Range("M111").Value = -343080.992
$ws.Range("H112").Value = 3069.8823
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").Value = ""
$ws.Range("H117").Value = 48846.08
$ws.Range("J117").Value = 48846.08
$ws.Range("L117").Value = 48846.08
$ws.Range("N117").Value = -58024.08
$ws.Range("H125").Value = 7410826.5
$ws.Range("I125").Value = 2192.4
$ws.Range("K125").Value = 19731.6
$ws.Range("M125").Value = -17271.6
$ws.Range("H128").Value = 89564.78
$ws.Range("J128").Value = 89564.78
$ws.Range("L128").Value = 89564.78
$ws.Range("N128").Value = -99524.78
$ws.Range("H138").Value = 6908.2246
$ws.Range("J138").Value = 7901.325
$ws.Range("L138").Value = 23703.975
$ws.Range("N138").Value = -33983.975
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1143.4445
$ws.Range("J2").Value = 701.8
$ws.Range("L2").Value = 701.8
$ws.Range("N2").Value = -927.8
$ws.Range("H32").Value = 3679.853
$ws.Range("I32").Value = 3679.853
$ws.Range("K32").Value = 3679.853
$ws.Range("M32").Value = -3392.853
$ws.Range("H45").Value = 2761.9333
$ws.Range("I45").Value = 1857.3636
$ws.Range("J45").Value = 5249.5
$ws.Range("K45").Value = 1857.3636
$ws.Range("L45").Value = 5249.5
$ws.Range("M45").Value = -1480.3636
$ws.Range("N45").Value = -6003.5
$ws.Range("H74").Value = 2068.9
$ws.Range("I74").Value = 1728.1765
$ws.Range("K74").Value = 1728.1765
$ws.Range("M74").Value = -854.1765
$ws.Range("H77").Value = 2068.9
$ws.Range("I77").Value = 1728.1765
$ws.Range("K77").Value = 8640.8825
$ws.Range("M77").Value = -4272.8825
$ws.Range("H116").Value = 1143.4445
$ws.Range("J116").Value = 701.8
$ws.Range("L116").Value = 701.8
$ws.Range("N116").Value = -5289.8
$ws.Range("H122").Value = 3527.7896
$ws.Range("I122").Value = 2137.75
$ws.Range("K122").Value = 6413.25
$ws.Range("M122").Value = -3963.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1143.4445
$ws.Range("J3").Value = 701.8
$ws.Range("L3").Value = 701.8
$ws.Range("N3").Value = -929.8
$ws.Range("H64").Value = 966.8570999999999
$ws.Range("J64").Value = 713.6
$ws.Range("L64").Value = 713.6
$ws.Range("N64").Value = -1163.6
$ws.Range("H67").Value = 966.8570999999999
$ws.Range("J67").Value = 713.6
$ws.Range("L67").Value = 713.6
$ws.Range("N67").Value = -2273.6
$ws.Range("H134").Value = 28057.924
$ws.Range("I134").Value = 1704.8182
$ws.Range("K134").Value = 5114.4546
$ws.Range("M134").Value = -2579.4546
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4132.968
$ws.Range("I132").Value = 3541.1904
$ws.Range("K132").Value = 10623.5712
$ws.Range("M132").Value = -8093.5712
$ws.Range("H141").Value = 223764.25
$ws.Range("J141").Value = 238545.42
$ws.Range("L141").Value = 238545.42
$ws.Range("N141").Value = -248905.42
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 500
$ws.Range("J7").Value = 500
$ws.Range("L7").Value = 1500
$ws.Range("N7").Value = -1724
$ws.Range("H13").Value = 2656.25
$ws.Range("J13").Value = 4212.5
$ws.Range("L13").Value = 12637.5
$ws.Range("N13").Value = -12973.5
$ws.Range("H41").Value = 1133.3334
$ws.Range("I41").Value = 1300
$ws.Range("J41").Value = 1050
$ws.Range("K41").Value = 3900
$ws.Range("L41").Value = 3150
$ws.Range("M41").Value = -3562
$ws.Range("N41").Value = -3826
$ws.Range("H125").Value = 3265
$ws.Range("I125").Value = 730
$ws.Range("K125").Value = 2190
$ws.Range("M125").Value = 2730
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1115358.8
$ws.Range("I80").Value = 773087.9399999999
$ws.Range("K80").Value = 773087.9399999999
$ws.Range("M80").Value = -772089.9399999999
$ws.Range("H83").Value = 1115358.8
$ws.Range("I83").Value = 773087.9399999999
$ws.Range("K83").Value = 3865439.7
$ws.Range("M83").Value = -3860447.7
$ws.Range("H102").Value = 1720.2
$ws.Range("I102").Value = 460.33334
$ws.Range("K102").Value = 460.33334
$ws.Range("M102").Value = 1161.66666
$ws.Range("H122").Value = 3102.8276
$ws.Range("J122").Value = 4478.125
$ws.Range("L122").Value = 13434.375
$ws.Range("N122").Value = -18334.375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 253829.75
$ws.Range("I7").Value = 3683.2
$ws.Range("K7").Value = 3683.2
$ws.Range("M7").Value = -3571.2
$ws.Range("H55").Value = 661.1818
$ws.Range("I55").Value = 186.96153
$ws.Range("K55").Value = 186.96153
$ws.Range("M55").Value = -13.96153000000001
$ws.Range("H61").Value = 4850.7144
$ws.Range("I61").Value = 4460.154
$ws.Range("J61").Value = 5485.375
$ws.Range("K61").Value = 4460.154
$ws.Range("L61").Value = 5485.375
$ws.Range("M61").Value = -4258.154
$ws.Range("N61").Value = -5889.375
$ws.Range("H82").Value = 650
$ws.Range("I82").Value = 650
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 650
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -289
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 650
$ws.Range("I85").Value = 650
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 650
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 598
$ws.Range("N85").Value = ""
$ws.Range("H113").Value = 4850.7144
$ws.Range("I113").Value = 4460.154
$ws.Range("J113").Value = 5485.375
$ws.Range("K113").Value = 4460.154
$ws.Range("L113").Value = 5485.375
$ws.Range("M113").Value = -2290.154
$ws.Range("N113").Value = -9825.375
$ws.Range("H122").Value = 502099.84
$ws.Range("I122").Value = 2063.3125
$ws.Range("J122").Value = 2502246
$ws.Range("K122").Value = 6189.9375
$ws.Range("L122").Value = 7506738
$ws.Range("M122").Value = -3739.9375
$ws.Range("N122").Value = -7511638
$ws.Range("H126").Value = 253829.75
$ws.Range("I126").Value = 3683.2
$ws.Range("K126").Value = 11049.6
$ws.Range("M126").Value = -8579.599999999999
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6823.909
$ws.Range("I74").Value = 4656.3335
$ws.Range("J74").Value = 7636.75
$ws.Range("K74").Value = 4656.3335
$ws.Range("L74").Value = 7636.75
$ws.Range("M74").Value = -3720.3335
$ws.Range("N74").Value = -9508.75
$ws.Range("H77").Value = 6823.909
$ws.Range("I77").Value = 4656.3335
$ws.Range("J77").Value = 7636.75
$ws.Range("K77").Value = 13969.0005
$ws.Range("L77").Value = 22910.25
$ws.Range("M77").Value = -9289.000499999998
$ws.Range("N77").Value = -32270.25
$ws.Range("H100").Value = 832.8
$ws.Range("I100").Value = 889.17645
$ws.Range("K100").Value = 1778.3529
$ws.Range("M100").Value = -1237.3529
$ws.Range("H122").Value = 34485416
$ws.Range("I122").Value = 40002140
$ws.Range("K122").Value = 120006420
$ws.Range("M122").Value = -120003970
$ws.Range("H124").Value = 84988.2
$ws.Range("J124").Value = 84988.2
$ws.Range("L124").Value = 84988.2
$ws.Range("N124").Value = -94808.2
$ws.Range("H126").Value = 5197.8
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5197.8
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15593.4
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -20533.4
$ws.Range("H136").Value = 12422817
$ws.Range("I136").Value = 20835568
$ws.Range("K136").Value = 62506704
$ws.Range("M136").Value = -62504154
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
